$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 8403740
$ws.Range("I2").Value = 9804030
$ws.Range("K2").Value = 9804030
$ws.Range("M2").Value = -9803917
$ws.Range("H111").Value = 6350.7144
$ws.Range("I111").Value = 5742.6665
$ws.Range("K111").Value = 17227.9995
$ws.Range("M111").Value = -14160.9995
$ws.Range("H116").Value = 3907.1177
$ws.Range("I116").Value = 3114
$ws.Range("K116").Value = 3114
$ws.Range("M116").Value = 328
$ws.Range("H141").Value = 16686
$ws.Range("I141").Value = 18784
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 56352
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -51172
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1859.2858
$ws.Range("I32").Value = 1766.9117
$ws.Range("K32").Value = 1766.9117
$ws.Range("M32").Value = -1479.9117
$ws.Range("H45").Value = 2117
$ws.Range("I45").Value = 2233.125
$ws.Range("K45").Value = 2233.125
$ws.Range("M45").Value = -1856.125
$ws.Range("H110").Value = 25552
$ws.Range("I110").Value = 25552
$ws.Range("K110").Value = 25552
$ws.Range("M110").Value = -23507
$ws.Range("H132").Value = 3955.5
$ws.Range("I132").Value = 3955.5
$ws.Range("K132").Value = 11866.5
$ws.Range("M132").Value = -9336.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 29000
$ws.Range("J56").Value = 29000
$ws.Range("L56").Value = 29000
$ws.Range("N56").Value = -30478
$ws.Range("H80").Value = 5226.067
$ws.Range("I80").Value = 334.85715
$ws.Range("K80").Value = 334.85715
$ws.Range("M80").Value = 663.14285
$ws.Range("H83").Value = 5226.067
$ws.Range("I83").Value = 334.85715
$ws.Range("K83").Value = 1674.28575
$ws.Range("M83").Value = 3317.71425
$ws.Range("H105").Value = 3109.8462
$ws.Range("I105").Value = 3142.9
$ws.Range("J105").Value = 2999.6667
$ws.Range("K105").Value = 3142.9
$ws.Range("L105").Value = 2999.6667
$ws.Range("M105").Value = -1395.9
$ws.Range("N105").Value = -6493.6667
$ws.Range("H134").Value = 9916.700000000001
$ws.Range("I134").Value = 10479.571
$ws.Range("J134").Value = 8603.333000000001
$ws.Range("K134").Value = 31438.713
$ws.Range("L134").Value = 25809.999
$ws.Range("M134").Value = -28903.713
$ws.Range("N134").Value = -30879.999
$ws.Range("H135").Value = 99996
$ws.Range("J135").Value = 99996
$ws.Range("L135").Value = 99996
$ws.Range("N135").Value = -110136

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19490
$ws.Range("I41").Value = 19490
$ws.Range("K41").Value = 19490
$ws.Range("M41").Value = -19062
$ws.Range("H132").Value = 2509.7273
$ws.Range("I132").Value = 2510.7
$ws.Range("K132").Value = 7532.099999999999
$ws.Range("M132").Value = -5002.099999999999
$ws.Range("H134").Value = 3445.5264
$ws.Range("I134").Value = 3360.3125
$ws.Range("K134").Value = 10080.9375
$ws.Range("M134").Value = -7545.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2233.7
$ws.Range("J129").Value = 3350
$ws.Range("L129").Value = 10050
$ws.Range("N129").Value = -20050
$ws.Range("H140").Value = 628436.4
$ws.Range("I140").Value = 772614
$ws.Range("K140").Value = 2317842
$ws.Range("M140").Value = -2312662
$ws.Range("H141").Value = 8145.7144
$ws.Range("I141").Value = 8145.7144
$ws.Range("K141").Value = 24437.1432
$ws.Range("M141").Value = -19257.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 216916.5
$ws.Range("I3").Value = 500450
$ws.Range("K3").Value = 500450
$ws.Range("M3").Value = -500334
$ws.Range("H14").Value = 1058027.2
$ws.Range("I14").Value = 1322500.2
$ws.Range("K14").Value = 1322500.2
$ws.Range("M14").Value = -1322332.2
$ws.Range("H97").Value = 726.3
$ws.Range("I97").Value = 701.3684
$ws.Range("K97").Value = 701.3684
$ws.Range("M97").Value = -205.3684
$ws.Range("H102").Value = 3416.4375
$ws.Range("I102").Value = 3377.6
$ws.Range("K102").Value = 3377.6
$ws.Range("M102").Value = -1755.6
$ws.Range("H132").Value = 2315.4614
$ws.Range("I132").Value = 2216.3333
$ws.Range("J132").Value = 2538.5
$ws.Range("K132").Value = 6648.999899999999
$ws.Range("L132").Value = 7615.5
$ws.Range("M132").Value = -4118.999899999999
$ws.Range("N132").Value = -12675.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3372.75
$ws.Range("I7").Value = 3121.8572
$ws.Range("J7").Value = 5129
$ws.Range("K7").Value = 3121.8572
$ws.Range("L7").Value = 5129
$ws.Range("M7").Value = -3009.8572
$ws.Range("N7").Value = -5353
$ws.Range("H40").Value = 1327.4117
$ws.Range("I40").Value = 1222.875
$ws.Range("K40").Value = 1222.875
$ws.Range("M40").Value = -1086.875
$ws.Range("H100").Value = 1852.2667
$ws.Range("I100").Value = 1746.7778
$ws.Range("K100").Value = 1746.7778
$ws.Range("M100").Value = -1205.7778
$ws.Range("H122").Value = 3159.7058
$ws.Range("I122").Value = 3169.6875
$ws.Range("K122").Value = 9509.0625
$ws.Range("M122").Value = -7059.0625
$ws.Range("H126").Value = 3372.75
$ws.Range("I126").Value = 3121.8572
$ws.Range("J126").Value = 5129
$ws.Range("K126").Value = 9365.571599999999
$ws.Range("L126").Value = 15387
$ws.Range("M126").Value = -6895.571599999999
$ws.Range("N126").Value = -20327
$ws.Range("H130").Value = 84950
$ws.Range("J130").Value = 84950
$ws.Range("L130").Value = 84950
$ws.Range("N130").Value = -94990
$ws.Range("H132").Value = 2122.4443
$ws.Range("I132").Value = 2085.077
$ws.Range("K132").Value = 6255.231000000001
$ws.Range("M132").Value = -3725.231000000001
$ws.Range("H136").Value = 5934.25
$ws.Range("I136").Value = 3271.1
$ws.Range("K136").Value = 9813.299999999999
$ws.Range("M136").Value = -7263.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 47887.69
$ws.Range("J45").Value = 61081.5
$ws.Range("L45").Value = 61081.5
$ws.Range("N45").Value = -62063.5
$ws.Range("H122").Value = 3499.3125
$ws.Range("I122").Value = 3619.9333
$ws.Range("K122").Value = 10859.7999
$ws.Range("M122").Value = -8409.7999
$ws.Range("H126").Value = 1170.5264
$ws.Range("I126").Value = 1116.5
$ws.Range("K126").Value = 3349.5
$ws.Range("M126").Value = -879.5
$ws.Range("H132").Value = 2202.5833
$ws.Range("I132").Value = 2129.6365
$ws.Range("K132").Value = 6388.9095
$ws.Range("M132").Value = -3858.9095
